# Update the "Generacion 1" results table (rows 2-51, columns B:E) with the
# refreshed genetic-algorithm output values. Source values are formatted
# numeric strings (e.g. "12,151.00"), so the destination range is forced to
# Text format before the assignment to keep Excel from auto-converting them
# to numbers (which would silently strip the thousands separators / fixed
# decimal formatting baked into the original text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, "12,151.00", "16,882.00", "15,055.38", "752,769.00"),
    @(3, "12,384.00", "16,882.00", "14,912.80", "745,640.00"),
    @(4, "12,148.00", "16,934.00", "14,870.38", "743,519.00"),
    @(5, "12,117.00", "16,934.00", "15,078.40", "753,920.00"),
    @(6, "13,257.00", "17,489.00", "15,529.62", "776,481.00"),
    @(7, "13,978.00", "18,454.00", "15,967.46", "798,373.00"),
    @(8, "14,128.00", "18,633.00", "16,282.16", "814,108.00"),
    @(9, "13,105.00", "18,747.00", "16,551.86", "827,593.00"),
    @(10, "13,818.00", "18,747.00", "16,474.62", "823,731.00"),
    @(11, "14,041.00", "18,747.00", "16,851.14", "842,557.00"),
    @(12, "14,606.00", "18,747.00", "16,831.32", "841,566.00"),
    @(13, "14,726.00", "18,963.00", "16,980.56", "849,028.00"),
    @(14, "14,920.00", "18,963.00", "16,858.60", "842,930.00"),
    @(15, "14,393.00", "18,963.00", "16,997.54", "849,877.00"),
    @(16, "14,632.00", "18,963.00", "17,103.64", "855,182.00"),
    @(17, "15,022.00", "19,133.00", "17,261.18", "863,059.00"),
    @(18, "15,092.00", "19,133.00", "17,233.20", "861,660.00"),
    @(19, "15,556.00", "19,133.00", "17,189.40", "859,470.00"),
    @(20, "15,265.00", "19,133.00", "17,226.14", "861,307.00"),
    @(21, "15,265.00", "19,133.00", "17,412.74", "870,637.00"),
    @(22, "16,019.00", "19,133.00", "17,524.00", "876,200.00"),
    @(23, "15,670.00", "19,391.00", "17,622.30", "881,115.00"),
    @(24, "15,851.00", "19,391.00", "17,528.14", "876,407.00"),
    @(25, "15,243.00", "19,391.00", "17,614.60", "880,730.00"),
    @(26, "15,367.00", "19,391.00", "17,682.76", "884,138.00"),
    @(27, "16,405.00", "19,396.00", "17,740.38", "887,019.00"),
    @(28, "15,975.00", "19,396.00", "17,896.24", "894,812.00"),
    @(29, "15,266.00", "19,396.00", "17,773.72", "888,686.00"),
    @(30, "16,070.00", "19,396.00", "17,743.22", "887,161.00"),
    @(31, "15,846.00", "19,396.00", "17,753.18", "887,659.00"),
    @(32, "16,341.00", "19,396.00", "17,778.60", "888,930.00"),
    @(33, "16,353.00", "19,396.00", "17,831.64", "891,582.00"),
    @(34, "16,601.00", "19,396.00", "18,003.22", "900,161.00"),
    @(35, "16,656.00", "19,396.00", "18,037.08", "901,854.00"),
    @(36, "16,696.00", "19,396.00", "18,061.48", "903,074.00"),
    @(37, "16,826.00", "19,396.00", "18,271.08", "913,554.00"),
    @(38, "16,635.00", "19,442.00", "18,519.64", "925,982.00"),
    @(39, "17,005.00", "19,442.00", "18,502.72", "925,136.00"),
    @(40, "17,051.00", "19,442.00", "18,445.68", "922,284.00"),
    @(41, "16,826.00", "19,442.00", "18,437.30", "921,865.00"),
    @(42, "16,780.00", "19,442.00", "18,572.74", "928,637.00"),
    @(43, "16,610.00", "19,442.00", "18,485.92", "924,296.00"),
    @(44, "16,494.00", "19,442.00", "18,520.18", "926,009.00"),
    @(45, "16,008.00", "19,442.00", "18,567.48", "928,374.00"),
    @(46, "16,301.00", "19,442.00", "18,479.66", "923,983.00"),
    @(47, "16,661.00", "19,442.00", "18,446.74", "922,337.00"),
    @(48, "16,179.00", "19,442.00", "18,423.80", "921,190.00"),
    @(49, "17,053.00", "19,442.00", "18,458.60", "922,930.00"),
    @(50, "15,960.00", "19,442.00", "18,423.38", "921,169.00"),
    @(51, "16,920.00", "19,667.00", "18,492.22", "924,611.00")
)

$targetRange = $ws.Range("B2:E51")
$targetRange.NumberFormat = "@"

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# The values were entered as literal text (they keep thousands separators and
# fixed decimal places, e.g. "752,769.00"), so the temporary Text number
# format is no longer needed once the literal strings are locked in. Clear it
# back off so the cells end up with their original (default) formatting,
# matching the source data which carried no explicit cell style.
$targetRange.ClearFormats()

Write-Host "Updated $($data.Count) rows in 'Generacion 1'"
